$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ACHData")

# --- Add new ACH test data rows (6, 7, 8) ---
# Row 6: ID=5, AccountNumber/Confirm=1234567, RoutingNumber=256072691, Notes=PCPM
$ws.Range("A6").Value = "5"
$ws.Range("B6").Value = "1234567"
$ws.Range("G6").Value = "PCPM"
$ws.Range("C6").Value = "1234567"
$ws.Range("D6").Value = "256072691"

# Row 7: ID=6, AccountNumber/Confirm=12345, RoutingNumber=256072691, Notes=PSPM
$ws.Range("A7").Value = "6"
$ws.Range("G7").Value = "PSPM"
$ws.Range("B7").Value = "12345"
$ws.Range("C7").Value = "12345"
$ws.Range("D7").Value = "256072691"

# Row 8: ID=7, AccountNumber/Confirm=123456, RoutingNumber=256072691, Notes=CorpPM
$ws.Range("B8").Value = "123456"
$ws.Range("A8").Value = "7"
$ws.Range("G8").Value = "CorpPM"
$ws.Range("C8").Value = "123456"
$ws.Range("D8").Value = "256072691"

# --- Update the active sheet / selection state ---
# Previously NameData (sheet1) was the selected tab; now ACHData (sheet3) is,
# with a new active selection cell.
$ws.Activate()
$ws.Range("E14").Select()
